$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

[object[,]]$data = New-Object "object[,]" 4,23

$data[0,0] = 0.00666173205033309
$data[0,1] = 0.0111028867505551
$data[0,2] = 0.00296076980014804
$data[0,3] = 0.00148038490007402
$data[0,4] = 0.0111028867505551
$data[0,5] = 0.117690599555885
$data[0,6] = 0.0399703923019985
$data[0,7] = 0.0162842339008142
$data[0,8] = 0.00370096225018505
$data[0,9] = 0.0125832716506292
$data[0,10] = 0.00518134715025907
$data[0,11] = 0.00370096225018505
$data[0,12] = 0.00370096225018505
$data[0,13] = 0.0185048112509252
$data[0,14] = 0.00148038490007402
$data[0,15] = 0.00666173205033309
$data[0,16] = 0.977054034048853
$data[0,17] = 0.0074019245003701
$data[0,18] = 0.0429311621021466
$data[0,19] = 0.0266469282013323
$data[0,20] = 0.00814211695040711
$data[0,21] = 0.00148038490007402
$data[0,22] = 0.0125832716506292
$data[1,0] = 0.00222057735011103
$data[1,1] = 0.686158401184308
$data[1,2] = 0.00222057735011103
$data[1,3] = 0.0421909696521095
$data[1,4] = 0.00074019245003701
$data[1,5] = 0.00888230940044412
$data[1,6] = 0.920059215396003
$data[1,7] = 0.973353071798668
$data[1,8] = 0.00074019245003701
$data[1,9] = 0.00074019245003701
$data[1,10] = 0.00666173205033309
$data[1,11] = 0.00592153960029608
$data[1,12] = 0.974833456698742
$data[1,13] = 0.00074019245003701
$data[1,14] = 0
$data[1,15] = 0
$data[1,16] = 0.00148038490007402
$data[1,17] = 0.982975573649149
$data[1,18] = 0.0895632864544782
$data[1,19] = 0.00518134715025907
$data[1,20] = 0
$data[1,21] = 0.00962250185048113
$data[1,22] = 0.00370096225018505
$data[2,0] = 0.990377498149519
$data[2,1] = 0.00370096225018505
$data[2,2] = 0.988156920799408
$data[2,3] = 0.0133234641006662
$data[2,4] = 0.982235381199112
$data[2,5] = 0.86380458919319
$data[2,6] = 0.00296076980014804
$data[2,7] = 0.00074019245003701
$data[2,8] = 0.0259067357512953
$data[2,9] = 0.981495188749075
$data[2,10] = 0.987416728349371
$data[2,11] = 0.988156920799408
$data[2,12] = 0.00222057735011103
$data[2,13] = 0.978534418948927
$data[2,14] = 0.998519615099926
$data[2,15] = 0.99259807549963
$data[2,16] = 0.0162842339008142
$data[2,17] = 0.00074019245003701
$data[2,18] = 0.0259067357512953
$data[2,19] = 0.962250185048113
$data[2,20] = 0.983715766099186
$data[2,21] = 0.987416728349371
$data[2,22] = 0.97779422649889
$data[3,0] = 0.00074019245003701
$data[3,1] = 0.299037749814952
$data[3,2] = 0.00666173205033309
$data[3,3] = 0.94300518134715
$data[3,4] = 0.00592153960029608
$data[3,5] = 0.00888230940044412
$data[3,6] = 0.0362694300518135
$data[3,7] = 0.00962250185048113
$data[3,8] = 0.969652109548483
$data[3,9] = 0.00444115470022206
$data[3,10] = 0.00074019245003701
$data[3,11] = 0.00074019245003701
$data[3,12] = 0.0192450037009623
$data[3,13] = 0.00148038490007402
$data[3,14] = 0
$data[3,15] = 0.00074019245003701
$data[3,16] = 0.00518134715025907
$data[3,17] = 0.00888230940044412
$data[3,18] = 0.84159881569208
$data[3,19] = 0.00592153960029608
$data[3,20] = 0.00814211695040711
$data[3,21] = 0
$data[3,22] = 0.00444115470022206

$ws.Range("B2:X5").Value = $data
